# Apply updated crypto market data to Sheet1 (A1:E51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.018.27'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''319.29'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '''100.61'
$ws.Range('E6').Value = '  +1.71%  '
$ws.Range('D7').Value = '''0.575'
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -2.59%  '
$ws.Range('D10').Value = '''36.87'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('D11').Value = '''0.0827'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('E12').Value = '  -2.42%  '
$ws.Range('E13').Value = '  -1.67%  '
$ws.Range('D14').Value = '2.588.41'
$ws.Range('E14').Value = '  +0.94%  '
$ws.Range('D15').Value = '''0.852'
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('D16').Value = '''14.30'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '2.248.18'
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').Value = '43.917.75'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('D19').Value = '''13.35'
$ws.Range('E19').Value = '  -3.75%  '
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').Value = '''6.44'
$ws.Range('E21').Value = '  -1.84%  '
$ws.Range('D22').Value = '''65.39'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('D23').Value = '''3.09'
$ws.Range('E23').Value = '  -4.23%  '
$ws.Range('D24').Value = '''233.56'
$ws.Range('E24').Value = '  -1.45%  '
$ws.Range('E25').Value = '  -5.45%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').Value = '''10.56'
$ws.Range('E27').Value = '  +5.10%  '
$ws.Range('D28').Value = '''38.84'
$ws.Range('E28').Value = '  +5.95%  '
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').Value = '''6.11'
$ws.Range('E30').Value = '  -5.28%  '
$ws.Range('D31').Value = '''158.86'
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('D33').Value = '''0.0843'
$ws.Range('E33').Value = '  -2.53%  '
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').Value = '''0.113'
$ws.Range('E35').Value = '  +8.41%  '
$ws.Range('D36').Value = '''3.10'
$ws.Range('E36').Value = '  -6.13%  '
$ws.Range('E37').Value = '  +3.87%  '
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('D39').Value = '''16.32'
$ws.Range('E39').Value = '  +14.64%  '
$ws.Range('E40').Value = '  -0.84%  '
$ws.Range('D41').Value = '''4.16'
$ws.Range('E41').Value = '  -5.99%  '
$ws.Range('E42').Value = '  -1.48%  '
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').Value = '1.768.52'
$ws.Range('E44').Value = '  +0.47%  '
$ws.Range('D45').Value = '''74.50'
$ws.Range('E45').Value = '  +0.72%  '
$ws.Range('E46').Value = '  -3.65%  '
$ws.Range('D47').Value = '''5.17'
$ws.Range('E47').Value = '  -2.27%  '
$ws.Range('D48').Value = '''81.02'
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '''1.68'
$ws.Range('E49').Value = '  +2.53%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '''103.28'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').Value = '''57.28'
$ws.Range('E51').Value = '  -1.12%  '
